$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134, shifting existing rows 134-253 down to 135-254
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new price-report record
$ws.Range("A134").Value = 8
$ws.Range("B134").Value = "Terminal La Palmera de La Serena"
$ws.Range("C134").Value = "Coquimbo"
$ws.Range("D134").Value = 44669
$ws.Range("E134").Value = 4
$ws.Range("F134").Value = 100112003
$ws.Range("G134").Value = "Ajo"
$ws.Range("H134").Value = "Chino"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 520
$ws.Range("K134").Value = 18500
$ws.Range("L134").Value = 19000
$ws.Range("M134").Value = 18750
$ws.Range("N134").Value = "`$/caja 10 kilos"
$ws.Range("O134").Value = "China"
$ws.Range("P134").Value = 1875
$ws.Range("Q134").Value = 10
$ws.Range("R134").Value = "Hortaliza"
